$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30.87085333333333
$ws.Range("H2").Value = 92.61256
$ws.Range("I2").Value = 0.2985789950947061
$ws.Range("J2").Value = 0.2985789950947061
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1149353333333333
$ws.Range("N2").Value = 0.344806
$ws.Range("O2").Value = 0.05719122335670149
$ws.Range("P2").Value = 0.05719122335670149
$ws.Range("Q2").Value = 3.548151818151111
$ws.Range("R2").Value = 31.93336636336
$ws.Range("S2").Value = 0.01707609799808081
$ws.Range("T2").Value = 0.01707609799808082
$ws.Range("G3").Value = 30.87085333333333
$ws.Range("H3").Value = 92.61256
$ws.Range("I3").Value = 0.2985789950947061
$ws.Range("J3").Value = 0.2985789950947061
$ws.Range("M3").Value = 1.894732
$ws.Range("N3").Value = 5.684196
$ws.Range("O3").Value = 0.9428087766432985
$ws.Range("P3").Value = 0.9428087766432984
$ws.Range("Q3").Value = 58.49199367797333
$ws.Range("R3").Value = 526.42794310176
$ws.Range("S3").Value = 0.2815028970966253
$ws.Range("T3").Value = 0.2815028970966253
$ws.Range("G4").Value = 33.793597
$ws.Range("I4").Value = 0.3268474027571036
$ws.Range("J4").Value = 0.3268474027571037
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1149353333333333
$ws.Range("N4").Value = 0.344806
$ws.Range("O4").Value = 0.05719122335670149
$ws.Range("P4").Value = 0.05719122335670149
$ws.Range("Q4").Value = 3.884078335727333
$ws.Range("R4").Value = 34.956705021546
$ws.Range("S4").Value = 0.01869280281463928
$ws.Range("T4").Value = 0.01869280281463928
$ws.Range("G5").Value = 33.793597
$ws.Range("I5").Value = 0.3268474027571036
$ws.Range("J5").Value = 0.3268474027571037
$ws.Range("M5").Value = 1.894732
$ws.Range("N5").Value = 5.684196
$ws.Range("O5").Value = 0.9428087766432985
$ws.Range("P5").Value = 0.9428087766432984
$ws.Range("Q5").Value = 64.029809631004
$ws.Range("R5").Value = 576.268286679036
$ws.Range("S5").Value = 0.3081545999424644
$ws.Range("T5").Value = 0.3081545999424644
$ws.Range("G6").Value = 2.981185666666667
$ws.Range("H6").Value = 8.943557
$ws.Range("I6").Value = 0.02883365130639111
$ws.Range("J6").Value = 0.02883365130639111
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1149353333333333
$ws.Range("N6").Value = 0.344806
$ws.Range("O6").Value = 0.05719122335670149
$ws.Range("P6").Value = 0.05719122335670149
$ws.Range("Q6").Value = 0.3426435683268889
$ws.Range("R6").Value = 3.083792114942
$ws.Range("S6").Value = 0.001649031792053062
$ws.Range("T6").Value = 0.001649031792053061
$ws.Range("G7").Value = 2.981185666666667
$ws.Range("H7").Value = 8.943557
$ws.Range("I7").Value = 0.02883365130639111
$ws.Range("J7").Value = 0.02883365130639111
$ws.Range("M7").Value = 1.894732
$ws.Range("N7").Value = 5.684196
$ws.Range("O7").Value = 0.9428087766432985
$ws.Range("P7").Value = 0.9428087766432984
$ws.Range("Q7").Value = 5.648547880574667
$ws.Range("R7").Value = 50.836930925172
$ws.Range("S7").Value = 0.02718461951433804
$ws.Range("T7").Value = 0.02718461951433804
$ws.Range("G8").Value = 35.74694633333333
$ws.Range("H8").Value = 107.240839
$ws.Range("I8").Value = 0.3457399508417991
$ws.Range("J8").Value = 0.3457399508417991
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1149353333333333
$ws.Range("N8").Value = 0.344806
$ws.Range("O8").Value = 0.05719122335670149
$ws.Range("P8").Value = 0.05719122335670149
$ws.Range("Q8").Value = 4.108587192470445
$ws.Range("R8").Value = 36.977284732234
$ws.Range("S8").Value = 0.01977329075192833
$ws.Range("T8").Value = 0.01977329075192832
$ws.Range("G9").Value = 35.74694633333333
$ws.Range("H9").Value = 107.240839
$ws.Range("I9").Value = 0.3457399508417991
$ws.Range("J9").Value = 0.3457399508417991
$ws.Range("M9").Value = 1.894732
$ws.Range("N9").Value = 5.684196
$ws.Range("O9").Value = 0.9428087766432985
$ws.Range("P9").Value = 0.9428087766432984
$ws.Range("Q9").Value = 67.73088312004934
$ws.Range("R9").Value = 609.577948080444
$ws.Range("S9").Value = 0.3259666600898708
$ws.Range("T9").Value = 0.3259666600898707
